$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.023873
$ws.Cells.Item(2, 8).Value = 0.071619
$ws.Cells.Item(2, 9).Value = 0.02747901635872243
$ws.Cells.Item(2, 10).Value = 0.02747901635872243
$ws.Cells.Item(2, 13).Value = 159.4836373333333
$ws.Cells.Item(2, 14).Value = 478.450912
$ws.Cells.Item(2, 15).Value = 0.2983285084902258
$ws.Cells.Item(2, 16).Value = 0.2983285084902258
$ws.Cells.Item(2, 17).Value = 3.807352874058668
$ws.Cells.Item(2, 18).Value = 34.266175866528
$ws.Cells.Item(2, 19).Value = 0.008197773965076177
$ws.Cells.Item(2, 20).Value = 0.008197773965076177

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.023873
$ws.Cells.Item(3, 8).Value = 0.071619
$ws.Cells.Item(3, 9).Value = 0.02747901635872243
$ws.Cells.Item(3, 10).Value = 0.02747901635872243
$ws.Cells.Item(3, 15).Value = 0.3227862111630279
$ws.Cells.Item(3, 16).Value = 0.3227862111630279
$ws.Cells.Item(3, 17).Value = 4.119488998881
$ws.Cells.Item(3, 18).Value = 37.075400989929
$ws.Cells.Item(3, 19).Value = 0.008869847576918876
$ws.Cells.Item(3, 20).Value = 0.008869847576918876

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.023873
$ws.Cells.Item(4, 8).Value = 0.071619
$ws.Cells.Item(4, 9).Value = 0.02747901635872243
$ws.Cells.Item(4, 10).Value = 0.02747901635872243
$ws.Cells.Item(4, 13).Value = 74.38770566666666
$ws.Cells.Item(4, 14).Value = 223.163117
$ws.Cells.Item(4, 15).Value = 0.1391489036280481
$ws.Cells.Item(4, 16).Value = 0.1391489036280482
$ws.Cells.Item(4, 17).Value = 1.775857697380333
$ws.Cells.Item(4, 18).Value = 15.982719276423
$ws.Cells.Item(4, 19).Value = 0.003823674999093425
$ws.Cells.Item(4, 20).Value = 0.003823674999093426

# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.023873
$ws.Cells.Item(5, 8).Value = 0.071619
$ws.Cells.Item(5, 9).Value = 0.02747901635872243
$ws.Cells.Item(5, 10).Value = 0.02747901635872243
$ws.Cells.Item(5, 13).Value = 58.41461433333333
$ws.Cells.Item(5, 14).Value = 175.243843
$ws.Cells.Item(5, 15).Value = 0.1092697975759847
$ws.Cells.Item(5, 16).Value = 0.1092697975759848
$ws.Cells.Item(5, 17).Value = 1.394532087979667
$ws.Cells.Item(5, 18).Value = 12.550788791817
$ws.Cells.Item(5, 19).Value = 0.003002626555104773
$ws.Cells.Item(5, 20).Value = 0.003002626555104773

# Row 6
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.023873
$ws.Cells.Item(6, 8).Value = 0.071619
$ws.Cells.Item(6, 9).Value = 0.02747901635872243
$ws.Cells.Item(6, 10).Value = 0.02747901635872243
$ws.Cells.Item(6, 13).Value = 69.746216
$ws.Cells.Item(6, 14).Value = 209.238648
$ws.Cells.Item(6, 15).Value = 0.1304665791427133
$ws.Cells.Item(6, 16).Value = 0.1304665791427133
$ws.Cells.Item(6, 17).Value = 1.665051414568
$ws.Cells.Item(6, 18).Value = 14.985462731112
$ws.Cells.Item(6, 19).Value = 0.003585093262529173
$ws.Cells.Item(6, 20).Value = 0.003585093262529174

# Row 7
$ws.Cells.Item(7, 9).Value = 0.04107483513127341
$ws.Cells.Item(7, 10).Value = 0.04107483513127341
$ws.Cells.Item(7, 13).Value = 159.4836373333333
$ws.Cells.Item(7, 14).Value = 478.450912
$ws.Cells.Item(7, 15).Value = 0.2983285084902258
$ws.Cells.Item(7, 16).Value = 0.2983285084902258
$ws.Cells.Item(7, 17).Value = 5.691120437027555
$ws.Cells.Item(7, 18).Value = 51.220083933248
$ws.Cells.Item(7, 19).Value = 0.01225379430119472
$ws.Cells.Item(7, 20).Value = 0.01225379430119473

# Row 8
$ws.Cells.Item(8, 9).Value = 0.04107483513127341
$ws.Cells.Item(8, 10).Value = 0.04107483513127341
$ws.Cells.Item(8, 15).Value = 0.3227862111630279
$ws.Cells.Item(8, 16).Value = 0.3227862111630279
$ws.Cells.Item(8, 19).Value = 0.01325839040616978
$ws.Cells.Item(8, 20).Value = 0.01325839040616978

# Row 9
$ws.Cells.Item(9, 9).Value = 0.04107483513127341
$ws.Cells.Item(9, 10).Value = 0.04107483513127341
$ws.Cells.Item(9, 13).Value = 74.38770566666666
$ws.Cells.Item(9, 14).Value = 223.163117
$ws.Cells.Item(9, 15).Value = 0.1391489036280481
$ws.Cells.Item(9, 16).Value = 0.1391489036280482
$ws.Cells.Item(9, 17).Value = 2.654500480813111
$ws.Cells.Item(9, 18).Value = 23.890504327318
$ws.Cells.Item(9, 19).Value = 0.005715518275219529
$ws.Cells.Item(9, 20).Value = 0.005715518275219531

# Row 10
$ws.Cells.Item(10, 9).Value = 0.04107483513127341
$ws.Cells.Item(10, 10).Value = 0.04107483513127341
$ws.Cells.Item(10, 13).Value = 58.41461433333333
$ws.Cells.Item(10, 14).Value = 175.243843
$ws.Cells.Item(10, 15).Value = 0.1092697975759847
$ws.Cells.Item(10, 16).Value = 0.1092697975759848
$ws.Cells.Item(10, 17).Value = 2.084506040946889
$ws.Cells.Item(10, 18).Value = 18.760554368522
$ws.Cells.Item(10, 19).Value = 0.004488238920261192
$ws.Cells.Item(10, 20).Value = 0.004488238920261193

# Row 11
$ws.Cells.Item(11, 9).Value = 0.04107483513127341
$ws.Cells.Item(11, 10).Value = 0.04107483513127341
$ws.Cells.Item(11, 13).Value = 69.746216
$ws.Cells.Item(11, 14).Value = 209.238648
$ws.Cells.Item(11, 15).Value = 0.1304665791427133
$ws.Cells.Item(11, 16).Value = 0.1304665791427133
$ws.Cells.Item(11, 17).Value = 2.488870469221333
$ws.Cells.Item(11, 18).Value = 22.399834222992
$ws.Cells.Item(11, 19).Value = 0.005358893228428183
$ws.Cells.Item(11, 20).Value = 0.005358893228428185

# Row 12
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.011261
$ws.Cells.Item(12, 8).Value = 0.033783
$ws.Cells.Item(12, 9).Value = 0.01296197391260307
$ws.Cells.Item(12, 10).Value = 0.01296197391260308
$ws.Cells.Item(12, 13).Value = 159.4836373333333
$ws.Cells.Item(12, 14).Value = 478.450912
$ws.Cells.Item(12, 15).Value = 0.2983285084902258
$ws.Cells.Item(12, 16).Value = 0.2983285084902258
$ws.Cells.Item(12, 17).Value = 1.795945240010667
$ws.Cells.Item(12, 18).Value = 16.163507160096
$ws.Cells.Item(12, 19).Value = 0.003866926344436092
$ws.Cells.Item(12, 20).Value = 0.003866926344436092

# Row 13
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.011261
$ws.Cells.Item(13, 8).Value = 0.033783
$ws.Cells.Item(13, 9).Value = 0.01296197391260307
$ws.Cells.Item(13, 10).Value = 0.01296197391260308
$ws.Cells.Item(13, 15).Value = 0.3227862111630279
$ws.Cells.Item(13, 16).Value = 0.3227862111630279
$ws.Cells.Item(13, 17).Value = 1.943181234717
$ws.Cells.Item(13, 18).Value = 17.488631112453
$ws.Cells.Item(13, 19).Value = 0.004183946448443155
$ws.Cells.Item(13, 20).Value = 0.004183946448443156

# Row 14
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.011261
$ws.Cells.Item(14, 8).Value = 0.033783
$ws.Cells.Item(14, 9).Value = 0.01296197391260307
$ws.Cells.Item(14, 10).Value = 0.01296197391260308
$ws.Cells.Item(14, 13).Value = 74.38770566666666
$ws.Cells.Item(14, 14).Value = 223.163117
$ws.Cells.Item(14, 15).Value = 0.1391489036280481
$ws.Cells.Item(14, 16).Value = 0.1391489036280482
$ws.Cells.Item(14, 17).Value = 0.8376799535123333
$ws.Cells.Item(14, 18).Value = 7.539119581611
$ws.Cells.Item(14, 19).Value = 0.001803644458794079
$ws.Cells.Item(14, 20).Value = 0.00180364445879408

# Row 15
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.011261
$ws.Cells.Item(15, 8).Value = 0.033783
$ws.Cells.Item(15, 9).Value = 0.01296197391260307
$ws.Cells.Item(15, 10).Value = 0.01296197391260308
$ws.Cells.Item(15, 13).Value = 58.41461433333333
$ws.Cells.Item(15, 14).Value = 175.243843
$ws.Cells.Item(15, 15).Value = 0.1092697975759847
$ws.Cells.Item(15, 16).Value = 0.1092697975759848
$ws.Cells.Item(15, 17).Value = 0.6578069720076667
$ws.Cells.Item(15, 18).Value = 5.920262748069
$ws.Cells.Item(15, 19).Value = 0.001416352265615333
$ws.Cells.Item(15, 20).Value = 0.001416352265615333

# Row 16
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.011261
$ws.Cells.Item(16, 8).Value = 0.033783
$ws.Cells.Item(16, 9).Value = 0.01296197391260307
$ws.Cells.Item(16, 10).Value = 0.01296197391260308
$ws.Cells.Item(16, 13).Value = 69.746216
$ws.Cells.Item(16, 14).Value = 209.238648
$ws.Cells.Item(16, 15).Value = 0.1304665791427133
$ws.Cells.Item(16, 16).Value = 0.1304665791427133
$ws.Cells.Item(16, 17).Value = 0.7854121383760001
$ws.Cells.Item(16, 18).Value = 7.068709245384
$ws.Cells.Item(16, 19).Value = 0.001691104395314414
$ws.Cells.Item(16, 20).Value = 0.001691104395314415

# Row 17
$ws.Cells.Item(17, 7).Value = 0.7979533333333334
$ws.Cells.Item(17, 8).Value = 2.39386
$ws.Cells.Item(17, 9).Value = 0.9184841745974011
$ws.Cells.Item(17, 10).Value = 0.9184841745974011
$ws.Cells.Item(17, 13).Value = 159.4836373333333
$ws.Cells.Item(17, 14).Value = 478.450912
$ws.Cells.Item(17, 15).Value = 0.2983285084902258
$ws.Cells.Item(17, 16).Value = 0.2983285084902258
$ws.Cells.Item(17, 17).Value = 127.2605000222578
$ws.Cells.Item(17, 18).Value = 1145.34450020032
$ws.Cells.Item(17, 19).Value = 0.2740100138795188
$ws.Cells.Item(17, 20).Value = 0.2740100138795188

# Row 18
$ws.Cells.Item(18, 7).Value = 0.7979533333333334
$ws.Cells.Item(18, 8).Value = 2.39386
$ws.Cells.Item(18, 9).Value = 0.9184841745974011
$ws.Cells.Item(18, 10).Value = 0.9184841745974011
$ws.Cells.Item(18, 15).Value = 0.3227862111630279
$ws.Cells.Item(18, 16).Value = 0.3227862111630279
$ws.Cells.Item(18, 17).Value = 137.69362787614
$ws.Cells.Item(18, 18).Value = 1239.24265088526
$ws.Cells.Item(18, 19).Value = 0.2964740267314961
$ws.Cells.Item(18, 20).Value = 0.2964740267314961

# Row 19
$ws.Cells.Item(19, 7).Value = 0.7979533333333334
$ws.Cells.Item(19, 8).Value = 2.39386
$ws.Cells.Item(19, 9).Value = 0.9184841745974011
$ws.Cells.Item(19, 10).Value = 0.9184841745974011
$ws.Cells.Item(19, 13).Value = 74.38770566666666
$ws.Cells.Item(19, 14).Value = 223.163117
$ws.Cells.Item(19, 15).Value = 0.1391489036280481
$ws.Cells.Item(19, 16).Value = 0.1391489036280482
$ws.Cells.Item(19, 17).Value = 59.35791769573556
$ws.Cells.Item(19, 18).Value = 534.22125926162
$ws.Cells.Item(19, 19).Value = 0.1278060658949411
$ws.Cells.Item(19, 20).Value = 0.1278060658949411

# Row 20
$ws.Cells.Item(20, 7).Value = 0.7979533333333334
$ws.Cells.Item(20, 8).Value = 2.39386
$ws.Cells.Item(20, 9).Value = 0.9184841745974011
$ws.Cells.Item(20, 10).Value = 0.9184841745974011
$ws.Cells.Item(20, 13).Value = 58.41461433333333
$ws.Cells.Item(20, 14).Value = 175.243843
$ws.Cells.Item(20, 15).Value = 0.1092697975759847
$ws.Cells.Item(20, 16).Value = 0.1092697975759848
$ws.Cells.Item(20, 17).Value = 46.61213622266445
$ws.Cells.Item(20, 18).Value = 419.50922600398
$ws.Cells.Item(20, 19).Value = 0.1003625798350034
$ws.Cells.Item(20, 20).Value = 0.1003625798350035

# Row 21
$ws.Cells.Item(21, 7).Value = 0.7979533333333334
$ws.Cells.Item(21, 8).Value = 2.39386
$ws.Cells.Item(21, 9).Value = 0.9184841745974011
$ws.Cells.Item(21, 10).Value = 0.9184841745974011
$ws.Cells.Item(21, 13).Value = 69.746216
$ws.Cells.Item(21, 14).Value = 209.238648
$ws.Cells.Item(21, 15).Value = 0.1304665791427133
$ws.Cells.Item(21, 16).Value = 0.1304665791427133
$ws.Cells.Item(21, 17).Value = 55.65422554458667
$ws.Cells.Item(21, 18).Value = 500.88802990128
$ws.Cells.Item(21, 19).Value = 0.1198314882564415
$ws.Cells.Item(21, 20).Value = 0.1198314882564416
